# Apply the "output generated at 456a3b4" gh-pages data refresh to
# 上海-漫展信息.xlsx (Shanghai expo/show info workbook).
#
# Sheet 1 "展览" (exhibitions)      - F-column (想去人数) counter bumps only
# Sheet 2 "演出" (live shows)       - F-column bumps + one new row inserted
#                                     (a new show announced for 2024-06-29)
# Sheet 3 "本地生活" (local life)   - F-column bumps only
# Sheet 4 "全部类型" (all types)    - mirrors the F-column bumps from the
#                                     sheets above (aggregated view; the
#                                     brand-new show row is NOT mirrored
#                                     here, matching the source diff)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: 展览
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")

$sheet1Updates = @{
    "F2"  = 1661
    "F3"  = 9519
    "F5"  = 756
    "F7"  = 234
    "F11" = 1421
    "F13" = 61
    "F14" = 1489
    "F15" = 130
    "F16" = 306
    "F18" = 149
    "F20" = 397
    "F22" = 103
    "F24" = 1
    "F25" = 50
    "F28" = 263
    "F29" = 77
    "F31" = 641
    "F33" = 2
    "F34" = 171
    "F36" = 185
    "F37" = 327
    "F39" = 269
    "F40" = 626
    "F42" = 746
    "F43" = 324
    "F45" = 323
    "F47" = 324
    "F49" = 16
}
foreach ($addr in $sheet1Updates.Keys) {
    $ws1.Range($addr).Value = $sheet1Updates[$addr]
}

# ---------------------------------------------------------------------------
# Sheet 2: 演出
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")

$sheet2UpdatesBeforeInsert = @{
    "F4"  = 38
    "F8"  = 10
    "F11" = 697
    "F19" = 974
    "F20" = 34
    "F22" = 285
}
foreach ($addr in $sheet2UpdatesBeforeInsert.Keys) {
    $ws2.Range($addr).Value = $sheet2UpdatesBeforeInsert[$addr]
}

# A brand-new show (2024-06-29) was published ahead of the existing
# 2024-06-30 entry, so a row is inserted at row 25 and everything below
# shifts down by one (old row 40 becomes row 41; dimension A1:I40 -> A1:I41).
$ws2.Rows.Item(25).Insert()

# The insert leaves the new row with a slightly different style for column
# A (no border/centering) than the rest of the table; copy the real
# formatting down from the row above so A25:I25 matches its siblings.
$ws2.Range("A24:I24").Copy()
$ws2.Range("A25:I25").PasteSpecial(-4122)   # xlPasteFormats

$ws2.Range("A25").Value = 24

# Column B holds plain-text dates ("2024-06-29"); without forcing a text
# format Excel would silently convert the entry into a date serial number.
$ws2.Range("B25").NumberFormat = "@"
$ws2.Range("B25").Value = "2024-06-29"

$ws2.Range("C25").Value = "上海·《沐云华·次元狂想》经典动漫二次元ACG音乐会—琥珀琴师×Mona×云小鱼"
$ws2.Range("D25").Value = "丁香路425号(上海科技馆地铁站1号口步行460米) 上海东方艺术中心音乐厅"
$ws2.Range("E25").Value = "2024.06.29 19:30-06.29 21:30"
$ws2.Range("F25").Value = 0
$ws2.Range("G25").Value = 99
$ws2.Range("H25").Value = "https://show.bilibili.com/platform/detail.html?id=86546"
$ws2.Range("I25").Value = "//i1.hdslb.com/bfs/openplatform/202405/obpyJwk21716875708282.png"

# The show that used to sit at (pre-insert) row 38 — Marcin Patrzalek — now
# lives at row 39; besides the shift its own counter also ticked up.
$ws2.Range("F39").Value = 23

# ---------------------------------------------------------------------------
# Sheet 3: 本地生活
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("本地生活")

$sheet3Updates = @{
    "F4"  = 774
    "F5"  = 374
    "F7"  = 2396
    "F8"  = 3655
    "F11" = 100
    "F12" = 115
}
foreach ($addr in $sheet3Updates.Keys) {
    $ws3.Range($addr).Value = $sheet3Updates[$addr]
}

# ---------------------------------------------------------------------------
# Sheet 4: 全部类型 (aggregated view; unaffected by the sheet-2 row insert)
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")

$sheet4Updates = @{
    "F2"  = 1661
    "F3"  = 774
    "F4"  = 9519
    "F5"  = 374
    "F7"  = 3655
    "F9"  = 756
    "F10" = 100
    "F11" = 100
    "F13" = 234
    "F16" = 697
    "F17" = 1421
    "F19" = 115
    "F20" = 1489
    "F21" = 130
    "F22" = 306
    "F24" = 149
    "F26" = 103
    "F28" = 50
    "F29" = 974
    "F31" = 34
    "F32" = 263
    "F34" = 285
    "F35" = 77
    "F37" = 641
    "F38" = 171
    "F39" = 324
    "F40" = 324
    "F41" = 327
    "F44" = 626
    "F46" = 746
    "F47" = 324
    "F50" = 323
    "F51" = 324
    "F52" = 23
}
foreach ($addr in $sheet4Updates.Keys) {
    $ws4.Range($addr).Value = $sheet4Updates[$addr]
}
